# Validation: Updated validation scripts to update Dunes 200 and ATF tests
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row labels first, so shared-string order matches the target ---
$ws.Range("A16").Value = "Smoke Alarm Activation Time"
$ws.Range("A17").Value = "Sprinkler Activation Time"

# --- Header text updates (row 1) ---
$ws.Range("F1").Value = "RP, 1824 only"
$ws.Range("J1").Value = "RP all"
$ws.Range("B1").Value = "Original from NUREG-1934"

# --- New style for J1:L1 (numFmt 0.00, center/center) ---
$ws.Range("J1:L1").NumberFormat = "0.00"
$ws.Range("J1:L1").HorizontalAlignment = -4108
$ws.Range("J1:L1").VerticalAlignment = -4108

# --- Row 2 style updates ---
$ws.Range("J2").NumberFormat = "0.00"

$ws.Range("K2:L2").NumberFormat = "0.00"
$ws.Range("K2:L2").HorizontalAlignment = -4108
$ws.Range("K2:L2").VerticalAlignment = -4108

# --- Row 4 updated values ---
$ws.Range("J4").Value = 0.98
$ws.Range("K4").Formula = "=0.45/2"

# --- Row 5 updated values ---
$ws.Range("J5").Value = 1.1599999999999999
$ws.Range("K5").Formula = "=0.43/2"

# --- Row 7: clear J7:L7 entirely ---
$ws.Range("J7:L7").Clear()

# --- Row 10: clear J10:K10 entirely ---
$ws.Range("J10:K10").Clear()

# --- Row 11: clear J11:K11 entirely ---
$ws.Range("J11:K11").Clear()

# --- Row 14 style updates ---
$ws.Range("J14").NumberFormat = "0.00"
$ws.Range("J14").VerticalAlignment = -4107

$ws.Range("L14").NumberFormat = "0.00"
$ws.Range("L14").VerticalAlignment = -4107

# --- Row 15 updated value ---
$ws.Range("J15").Value = 1

# --- Row 16: new values (label already set above) ---
$ws.Range("J16").Value = 1.05
$ws.Range("K16").Formula = "=0.98/2"
$ws.Range("L16").Formula = "=0.33/2"
$ws.Range("J16:L16").NumberFormat = "0.00"

# --- Row 17: new row with values (label already set above) ---
$ws.Range("J17").Value = 0.84
$ws.Range("K17").Formula = "=0.52/2"
$ws.Range("L17").Formula = "=0.2/2"
$ws.Range("J17:L17").NumberFormat = "0.00"

# --- Selection ---
$ws.Range("B6").Select()
